$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap mismatched match rows (F:V) introduced by corrected source scrape ---
# Row 14 <- old row 15 (F:V); Row 15 <- old row 14 (F:V)
$ws.Cells.Item(14,6).Value = "Sturm Graz"
$ws.Cells.Item(14,7).Value = 0
$ws.Cells.Item(14,8).Value = "A. Klagenfurt"
$ws.Cells.Item(14,9).Value = 0
$ws.Cells.Item(14,10).Value = 1.29
$ws.Cells.Item(14,11).Value = "05/08/2023 19:42"
$ws.Cells.Item(14,12).Value = 1.4
$ws.Cells.Item(14,13).Value = "12/08/2023 16:52"
$ws.Cells.Item(14,14).Value = 6.49
$ws.Cells.Item(14,15).Value = "05/08/2023 19:42"
$ws.Cells.Item(14,16).Value = 5.64
$ws.Cells.Item(14,17).Value = "12/08/2023 16:58"
$ws.Cells.Item(14,18).Value = 8.029999999999999
$ws.Cells.Item(14,19).Value = "05/08/2023 19:42"
$ws.Cells.Item(14,20).Value = 7.07
$ws.Cells.Item(14,21).Value = "12/08/2023 16:58"
$ws.Cells.Item(14,22).Value = "https://www.betexplorer.com/football/austria/bundesliga/sturm-graz-a-klagenfurt/p6ZpX4hm/"
$ws.Cells.Item(15,6).Value = "Wolfsberger AC"
$ws.Cells.Item(15,7).Value = 1
$ws.Cells.Item(15,8).Value = "A. Lustenau"
$ws.Cells.Item(15,9).Value = 1
$ws.Cells.Item(15,10).Value = 1.75
$ws.Cells.Item(15,11).Value = "06/08/2023 17:12"
$ws.Cells.Item(15,12).Value = 1.81
$ws.Cells.Item(15,13).Value = "12/08/2023 16:54"
$ws.Cells.Item(15,14).Value = 4.06
$ws.Cells.Item(15,15).Value = "06/08/2023 17:12"
$ws.Cells.Item(15,16).Value = 4.11
$ws.Cells.Item(15,17).Value = "12/08/2023 16:59"
$ws.Cells.Item(15,18).Value = 4.47
$ws.Cells.Item(15,19).Value = "06/08/2023 17:12"
$ws.Cells.Item(15,20).Value = 4.26
$ws.Cells.Item(15,21).Value = "12/08/2023 16:59"
$ws.Cells.Item(15,22).Value = "https://www.betexplorer.com/football/austria/bundesliga/wolfsberger-ac-a-lustenau/SUOkWO7g/"

# Row 20 <- old row 21 (F:V); Row 21 <- old row 20 (F:V)
$ws.Cells.Item(20,6).Value = "A. Lustenau"
$ws.Cells.Item(20,7).Value = 0
$ws.Cells.Item(20,8).Value = "Sturm Graz"
$ws.Cells.Item(20,9).Value = 1
$ws.Cells.Item(20,10).Value = 5.85
$ws.Cells.Item(20,11).Value = "12/08/2023 17:12"
$ws.Cells.Item(20,12).Value = 8.08
$ws.Cells.Item(20,13).Value = "19/08/2023 16:58"
$ws.Cells.Item(20,14).Value = 4.67
$ws.Cells.Item(20,15).Value = "12/08/2023 17:12"
$ws.Cells.Item(20,16).Value = 5.37
$ws.Cells.Item(20,17).Value = "19/08/2023 16:58"
$ws.Cells.Item(20,18).Value = 1.54
$ws.Cells.Item(20,19).Value = "12/08/2023 17:12"
$ws.Cells.Item(20,20).Value = 1.38
$ws.Cells.Item(20,21).Value = "19/08/2023 16:58"
$ws.Cells.Item(20,22).Value = "https://www.betexplorer.com/football/austria/bundesliga/a-lustenau-sturm-graz/b9aKo47t/"
$ws.Cells.Item(21,6).Value = "Tirol"
$ws.Cells.Item(21,7).Value = 1
$ws.Cells.Item(21,8).Value = "LASK"
$ws.Cells.Item(21,9).Value = 1
$ws.Cells.Item(21,10).Value = 3.25
$ws.Cells.Item(21,11).Value = "13/08/2023 17:12"
$ws.Cells.Item(21,12).Value = 5.89
$ws.Cells.Item(21,13).Value = "19/08/2023 16:56"
$ws.Cells.Item(21,14).Value = 3.9
$ws.Cells.Item(21,15).Value = "13/08/2023 17:12"
$ws.Cells.Item(21,16).Value = 4.89
$ws.Cells.Item(21,17).Value = "19/08/2023 16:56"
$ws.Cells.Item(21,18).Value = 2.1
$ws.Cells.Item(21,19).Value = "13/08/2023 17:12"
$ws.Cells.Item(21,20).Value = 1.52
$ws.Cells.Item(21,21).Value = "19/08/2023 16:55"
$ws.Cells.Item(21,22).Value = "https://www.betexplorer.com/football/austria/bundesliga/tirol-lask-linz/8pKAR0MO/"

# Row 23 <- old row 24 (F:V); Row 24 <- old row 23 (F:V)
$ws.Cells.Item(23,6).Value = "BW Linz"
$ws.Cells.Item(23,7).Value = 0
$ws.Cells.Item(23,8).Value = "Rapid Vienna"
$ws.Cells.Item(23,9).Value = 5
$ws.Cells.Item(23,10).Value = 4.16
$ws.Cells.Item(23,11).Value = "13/08/2023 17:12"
$ws.Cells.Item(23,12).Value = 3.74
$ws.Cells.Item(23,13).Value = "20/08/2023 16:55"
$ws.Cells.Item(23,14).Value = 3.95
$ws.Cells.Item(23,15).Value = "13/08/2023 17:12"
$ws.Cells.Item(23,16).Value = 3.95
$ws.Cells.Item(23,17).Value = "20/08/2023 16:55"
$ws.Cells.Item(23,18).Value = 1.85
$ws.Cells.Item(23,19).Value = "13/08/2023 17:12"
$ws.Cells.Item(23,20).Value = 1.97
$ws.Cells.Item(23,21).Value = "20/08/2023 16:55"
$ws.Cells.Item(23,22).Value = "https://www.betexplorer.com/football/austria/bundesliga/bw-linz-rapid-vienna/KOmTqrxg/"
$ws.Cells.Item(24,6).Value = "Hartberg"
$ws.Cells.Item(24,7).Value = 1
$ws.Cells.Item(24,8).Value = "Salzburg"
$ws.Cells.Item(24,9).Value = 5
$ws.Cells.Item(24,10).Value = 5.29
$ws.Cells.Item(24,11).Value = "13/08/2023 17:12"
$ws.Cells.Item(24,12).Value = 5.02
$ws.Cells.Item(24,13).Value = "20/08/2023 16:59"
$ws.Cells.Item(24,14).Value = 4.8
$ws.Cells.Item(24,15).Value = "13/08/2023 17:12"
$ws.Cells.Item(24,16).Value = 4.63
$ws.Cells.Item(24,17).Value = "20/08/2023 16:59"
$ws.Cells.Item(24,18).Value = 1.55
$ws.Cells.Item(24,19).Value = "13/08/2023 17:12"
$ws.Cells.Item(24,20).Value = 1.63
$ws.Cells.Item(24,21).Value = "20/08/2023 16:58"
$ws.Cells.Item(24,22).Value = "https://www.betexplorer.com/football/austria/bundesliga/hartberg-salzburg/Ea2Wr2ia/"


# --- Append new rows 44-49 ---
# Row 44
$ws.Cells.Item(44,1).Value = 43
$ws.Cells.Item(44,2).Value = "austria"
$ws.Cells.Item(44,3).Value = "bundesliga"
$ws.Cells.Item(44,4).Value = "2023-2024"
$ws.Cells.Item(44,5).Value = 45192.70833333334
$ws.Cells.Item(44,6).Value = "A. Klagenfurt"
$ws.Cells.Item(44,7).Value = 1
$ws.Cells.Item(44,8).Value = "A. Lustenau"
$ws.Cells.Item(44,9).Value = 0
$ws.Cells.Item(44,10).Value = 1.68
$ws.Cells.Item(44,11).Value = "17/09/2023 13:43"
$ws.Cells.Item(44,12).Value = 1.69
$ws.Cells.Item(44,13).Value = "23/09/2023 16:50"
$ws.Cells.Item(44,14).Value = 4.12
$ws.Cells.Item(44,15).Value = "17/09/2023 13:43"
$ws.Cells.Item(44,16).Value = 4.04
$ws.Cells.Item(44,17).Value = "23/09/2023 16:50"
$ws.Cells.Item(44,18).Value = 4.84
$ws.Cells.Item(44,19).Value = "17/09/2023 13:43"
$ws.Cells.Item(44,20).Value = 5.17
$ws.Cells.Item(44,21).Value = "23/09/2023 16:50"
$ws.Cells.Item(44,22).Value = "https://www.betexplorer.com/football/austria/bundesliga/a-klagenfurt-a-lustenau/Eu6h7xzH/"

# Row 45
$ws.Cells.Item(45,1).Value = 44
$ws.Cells.Item(45,2).Value = "austria"
$ws.Cells.Item(45,3).Value = "bundesliga"
$ws.Cells.Item(45,4).Value = "2023-2024"
$ws.Cells.Item(45,5).Value = 45192.70833333334
$ws.Cells.Item(45,6).Value = "Salzburg"
$ws.Cells.Item(45,7).Value = 0
$ws.Cells.Item(45,8).Value = "BW Linz"
$ws.Cells.Item(45,9).Value = 1
$ws.Cells.Item(45,10).Value = 1.16
$ws.Cells.Item(45,11).Value = "16/09/2023 18:42"
$ws.Cells.Item(45,12).Value = 1.26
$ws.Cells.Item(45,13).Value = "23/09/2023 16:50"
$ws.Cells.Item(45,14).Value = 8.359999999999999
$ws.Cells.Item(45,15).Value = "16/09/2023 18:42"
$ws.Cells.Item(45,16).Value = 6.59
$ws.Cells.Item(45,17).Value = "23/09/2023 16:57"
$ws.Cells.Item(45,18).Value = 14.49
$ws.Cells.Item(45,19).Value = "16/09/2023 18:42"
$ws.Cells.Item(45,20).Value = 10.89
$ws.Cells.Item(45,21).Value = "23/09/2023 16:57"
$ws.Cells.Item(45,22).Value = "https://www.betexplorer.com/football/austria/bundesliga/salzburg-bw-linz/6F2p9b54/"

# Row 46
$ws.Cells.Item(46,1).Value = 45
$ws.Cells.Item(46,2).Value = "austria"
$ws.Cells.Item(46,3).Value = "bundesliga"
$ws.Cells.Item(46,4).Value = "2023-2024"
$ws.Cells.Item(46,5).Value = 45192.70833333334
$ws.Cells.Item(46,6).Value = "Tirol"
$ws.Cells.Item(46,7).Value = 2
$ws.Cells.Item(46,8).Value = "Wolfsberger AC"
$ws.Cells.Item(46,9).Value = 3
$ws.Cells.Item(46,10).Value = 2.88
$ws.Cells.Item(46,11).Value = "17/09/2023 16:13"
$ws.Cells.Item(46,12).Value = 3.28
$ws.Cells.Item(46,13).Value = "23/09/2023 16:50"
$ws.Cells.Item(46,14).Value = 3.62
$ws.Cells.Item(46,15).Value = "17/09/2023 16:13"
$ws.Cells.Item(46,16).Value = 3.62
$ws.Cells.Item(46,17).Value = "23/09/2023 16:50"
$ws.Cells.Item(46,18).Value = 2.41
$ws.Cells.Item(46,19).Value = "17/09/2023 16:13"
$ws.Cells.Item(46,20).Value = 2.24
$ws.Cells.Item(46,21).Value = "23/09/2023 16:38"
$ws.Cells.Item(46,22).Value = "https://www.betexplorer.com/football/austria/bundesliga/tirol-wolfsberger-ac/SnnT0zDp/"

# Row 47
$ws.Cells.Item(47,1).Value = 46
$ws.Cells.Item(47,2).Value = "austria"
$ws.Cells.Item(47,3).Value = "bundesliga"
$ws.Cells.Item(47,4).Value = "2023-2024"
$ws.Cells.Item(47,5).Value = 45193.60416666666
$ws.Cells.Item(47,6).Value = "Altach"
$ws.Cells.Item(47,7).Value = 2
$ws.Cells.Item(47,8).Value = "Austria Vienna"
$ws.Cells.Item(47,9).Value = 1
$ws.Cells.Item(47,10).Value = 3.58
$ws.Cells.Item(47,11).Value = "17/09/2023 13:43"
$ws.Cells.Item(47,12).Value = 2.95
$ws.Cells.Item(47,13).Value = "24/09/2023 14:29"
$ws.Cells.Item(47,14).Value = 3.87
$ws.Cells.Item(47,15).Value = "17/09/2023 13:43"
$ws.Cells.Item(47,16).Value = 3.53
$ws.Cells.Item(47,17).Value = "24/09/2023 14:29"
$ws.Cells.Item(47,18).Value = 1.99
$ws.Cells.Item(47,19).Value = "17/09/2023 13:43"
$ws.Cells.Item(47,20).Value = 2.47
$ws.Cells.Item(47,21).Value = "24/09/2023 14:29"
$ws.Cells.Item(47,22).Value = "https://www.betexplorer.com/football/austria/bundesliga/altach-austria-vienna/zef15G4T/"

# Row 48
$ws.Cells.Item(48,1).Value = 47
$ws.Cells.Item(48,2).Value = "austria"
$ws.Cells.Item(48,3).Value = "bundesliga"
$ws.Cells.Item(48,4).Value = "2023-2024"
$ws.Cells.Item(48,5).Value = 45193.60416666666
$ws.Cells.Item(48,6).Value = "LASK"
$ws.Cells.Item(48,7).Value = 0
$ws.Cells.Item(48,8).Value = "Hartberg"
$ws.Cells.Item(48,9).Value = 0
$ws.Cells.Item(48,10).Value = 1.59
$ws.Cells.Item(48,11).Value = "17/09/2023 13:43"
$ws.Cells.Item(48,12).Value = 1.74
$ws.Cells.Item(48,13).Value = "24/09/2023 14:29"
$ws.Cells.Item(48,14).Value = 4.48
$ws.Cells.Item(48,15).Value = "17/09/2023 13:43"
$ws.Cells.Item(48,16).Value = 4.17
$ws.Cells.Item(48,17).Value = "24/09/2023 14:29"
$ws.Cells.Item(48,18).Value = 5.48
$ws.Cells.Item(48,19).Value = "17/09/2023 13:43"
$ws.Cells.Item(48,20).Value = 4.62
$ws.Cells.Item(48,21).Value = "24/09/2023 14:14"
$ws.Cells.Item(48,22).Value = "https://www.betexplorer.com/football/austria/bundesliga/lask-linz-hartberg/0x2l8IKA/"

# Row 49
$ws.Cells.Item(49,1).Value = 48
$ws.Cells.Item(49,2).Value = "austria"
$ws.Cells.Item(49,3).Value = "bundesliga"
$ws.Cells.Item(49,4).Value = "2023-2024"
$ws.Cells.Item(49,5).Value = 45193.70833333334
$ws.Cells.Item(49,6).Value = "Rapid Vienna"
$ws.Cells.Item(49,7).Value = 1
$ws.Cells.Item(49,8).Value = "Sturm Graz"
$ws.Cells.Item(49,9).Value = 1
$ws.Cells.Item(49,10).Value = 2.96
$ws.Cells.Item(49,11).Value = "17/09/2023 16:13"
$ws.Cells.Item(49,12).Value = 2.96
$ws.Cells.Item(49,13).Value = "24/09/2023 16:57"
$ws.Cells.Item(49,14).Value = 3.61
$ws.Cells.Item(49,15).Value = "17/09/2023 16:13"
$ws.Cells.Item(49,16).Value = 3.67
$ws.Cells.Item(49,17).Value = "24/09/2023 16:57"
$ws.Cells.Item(49,18).Value = 2.36
$ws.Cells.Item(49,19).Value = "17/09/2023 16:13"
$ws.Cells.Item(49,20).Value = 2.39
$ws.Cells.Item(49,21).Value = "24/09/2023 16:57"
$ws.Cells.Item(49,22).Value = "https://www.betexplorer.com/football/austria/bundesliga/rapid-vienna-sturm-graz/d4ec6dkN/"



# --- Apply header-matching style (bold, thin border, centered) to new index
#     column cells and the datetime number format to the new date column
#     cells, mirroring the formatting already used throughout column A / E. ---
for ($r = 44; $r -le 49; $r++) {
    $idxCell = $ws.Cells.Item($r, 1)
    $idxCell.Font.Name = "Calibri"
    $idxCell.Font.Size = 11
    $idxCell.Font.Bold = $true
    $idxCell.Borders.LineStyle = 1
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160

    $dateCell = $ws.Cells.Item($r, 5)
    $dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
